# Bulgaria First League workbook update (11-04-2024)
# - Swap data for row pair (8,9) and row pair (10,12)  [ids 6/7 and 8/10 got re-sorted]
# - Drop the stale/incomplete match that was on row 237 (CSKA 1948 Sofia vs Ludogorets,
#   B=7956727) and shift the remaining rows (238..243) up by one so the table ends at row 242
#   instead of row 243.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold actual match data (everything to the right of the sequential id in col A).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value()
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $data[$c]
    }
}

# --- Swap rows 8 and 9 ---
$row8 = Get-RowData 8
$row9 = Get-RowData 9
Set-RowData 8 $row9
Set-RowData 9 $row8

# --- Swap rows 10 and 12 ---
$row10 = Get-RowData 10
$row12 = Get-RowData 12
Set-RowData 10 $row12
Set-RowData 12 $row10

# --- Shift rows 238..243 up into 237..242 (row 237's original match is being removed) ---
for ($r = 237; $r -le 242; $r++) {
    $src = Get-RowData ($r + 1)
    Set-RowData $r $src
}

# --- Remove the now-duplicated last row ---
$ws.Rows.Item(243).Delete()
